# The commit re-syncs a set of match rows whose "id" (column B) and all of
# their associated stats (columns E..AB) had been attributed to the wrong
# row of an adjacent pair. Columns A (row #), C (Div) and D (Date) are
# identical for each pair and stay untouched; every other column between
# B and AB is swapped between the two rows listed below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")
$pairs = @(
    @(28,29),
    @(101,102),
    @(109,110),
    @(149,150),
    @(213,214),
    @(215,216),
    @(263,265),
    @(271,272),
    @(307,308)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}
